$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added for "Feria Lagunitas de Puerto Montt -
# Ciboulette" at row 121, pushing every following record down by one row
# (old row 121 -> new row 122, ..., old row 245 -> new row 246).
$ws.Rows.Item(121).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A121").Value = 4
$ws.Range("B121").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C121").Value = "Los Lagos"
$ws.Range("D121").Value = 44781
$ws.Range("E121").Value = 10
$ws.Range("F121").Value = 100112039
$ws.Range("G121").Value = "Ciboulette"
$ws.Range("H121").Value = "Sin especificar"
$ws.Range("I121").Value = "Primera"
$ws.Range("J121").Value = 80
$ws.Range("K121").Value = 4000
$ws.Range("L121").Value = 4000
$ws.Range("M121").Value = 4000
$ws.Range("N121").Value = "$/docena de atados"
$ws.Range("O121").Value = "Región Metropolitana"
$ws.Range("P121").Value = 1333
$ws.Range("Q121").Value = 3
$ws.Range("R121").Value = "Hortaliza"
